# Remove the "Date Placeholder" (dt, idx=10) shapes that were added to
# several slides. These are plain placeholder shapes showing the
# "MM.DD.20XX" boilerplate text; the revision removes them entirely.
#
# Note: on this host, Shape.Delete() on a <p:sp> shape does not actually
# remove the shape node from the tree (it renames/clears it instead and
# leaves the count unchanged) - Shape.Cut() performs a real removal, so
# we use that instead.

$p = $ppt.ActivePresentation

$targets = @(
    @{ Slide = 2; Name = "Date Placeholder 5" },
    @{ Slide = 3; Name = "Date Placeholder 1" },
    @{ Slide = 4; Name = "Date Placeholder 1" },
    @{ Slide = 6; Name = "Date Placeholder 3" },
    @{ Slide = 7; Name = "Date Placeholder 3" },
    @{ Slide = 8; Name = "Date Placeholder 1" },
    @{ Slide = 9; Name = "Date Placeholder 2" }
)

foreach ($t in $targets) {
    $slide = $p.Slides.Item($t.Slide)
    $shape = $slide.Shapes.Item($t.Name)
    $shape.Cut()
}
